$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$ws.Range("D2").Value = '72.556.20'
$ws.Range("E2").Value = '  +4.54%  '
$ws.Range("D3").Value = '4.070.60'
$ws.Range("E3").Value = '  +4.05%  '
$ws.Range("E4").Value = '  +0.20%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '520.57'
$ws.Range("E5").Value = '  -1.99%  '
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '148.01'
$ws.Range("E6").Value = '  +2.71%  '
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '0.738'
$ws.Range("E7").Value = '  +20.23%  '
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '1.00'
$ws.Range("E8").Value = '  +0.15%  '
$ws.Range("E9").Value = '  +8.19%  '
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '0.175'
$ws.Range("E10").Value = '  +1.22%  '
$ws.Range("E11").Value = '  -2.35%  '
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '47.61'
$ws.Range("E12").Value = '  +13.19%  '
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '11.22'
$ws.Range("E13").Value = '  +9.29%  '
$ws.Range("D14").Value = '4.721.92'
$ws.Range("E14").Value = '  +4.11%  '
$ws.Range("D15").Value = '4.069.78'
$ws.Range("E15").Value = '  +4.23%  '
$ws.Range("E16").Value = '  +7.86%  '
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '14.23'
$ws.Range("E17").Value = '  +1.72%  '
$ws.Range("E18").Value = '  -1.07%  '
$ws.Range("E19").Value = '  -1.43%  '
$ws.Range("D20").Value = '72.526.77'
$ws.Range("E20").Value = '  +4.66%  '
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '446.77'
$ws.Range("E21").Value = '  +4.35%  '
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '104.96'
$ws.Range("E22").Value = '  +18.36%  '
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '3.61'
$ws.Range("E23").Value = '  +6.71%  '
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '14.83'
$ws.Range("E24").Value = '  +4.93%  '
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '4.00'
$ws.Range("E25").Value = '  -1.71%  '
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '11.48'
$ws.Range("E26").Value = '  -0.05%  '
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '11.09'
$ws.Range("E27").Value = '  +4.61%  '
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '38.09'
$ws.Range("E28").Value = '  +4.38%  '
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '5.80'
$ws.Range("E29").Value = '  +2.06%  '
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '3.20'
$ws.Range("E30").Value = '  +13.29%  '
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '13.75'
$ws.Range("E31").Value = '  +4.67%  '
$ws.Range("E32").Value = '  +3.93%  '
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '680.16'
$ws.Range("E33").Value = '  +1.31%  '
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '6.86'
$ws.Range("E34").Value = '  +15.40%  '
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '67.31'
$ws.Range("E35").Value = '  -1.98%  '
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '43.79'
$ws.Range("E36").Value = '  +9.27%  '
$ws.Range("E37").Value = '  -1.95%  '
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '0.431'
$ws.Range("E38").Value = '  -1.37%  '
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '3.60'
$ws.Range("E39").Value = '  +11.87%  '
$ws.Range("E40").Value = '  +2.31%  '
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '1.00'
$ws.Range("E41").Value = '  +0.13%  '
$ws.Range("E42").Value = '  +3.89%  '
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '0.999'
$ws.Range("E43").Value = '  -0.10%  '
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '0.160'
$ws.Range("E44").Value = '  +13.81%  '
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '3.22'
$ws.Range("E45").Value = '  +1.11%  '
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '2.72'
$ws.Range("E46").Value = '  -3.19%  '
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '3.44'
$ws.Range("E47").Value = '  +0.71%  '
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '3.07'
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '9.18'
$ws.Range("E49").Value = '  +8.53%  '
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '3.36'
$ws.Range("E50").Value = '  +3.06%  '
$ws.Range("E51").Value = '  +1.51%  '
